$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns (and any text-like values) are stored as text,
# matching the source workbook where these cells are inline strings,
# not numbers/percentages - prevents Excel auto-converting "1.00" -> 1, etc.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.928.34'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.595.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.42'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.607.25'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.74'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.105'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.162'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.357'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.045.22'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.885.95'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.37'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000138'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.588.47'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.59'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.23'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.43'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.53'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.480'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.23'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.160'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.50'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0779'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.25'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.98%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.69'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '158.38'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.22'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.16'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.19'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.918'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.63'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.853'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.70'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '292.18'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.41'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0975'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.600'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0536'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.34%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.65'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0236'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.78'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.30%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.986.60'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.28%  '
